$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update column F (想去人数 / "want to go" count)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 265
$wsExhibit.Range("F5").Value = 234
$wsExhibit.Range("F6").Value = 266
$wsExhibit.Range("F7").Value = 67
$wsExhibit.Range("F10").Value = 47
$wsExhibit.Range("F11").Value = 31
$wsExhibit.Range("F13").Value = 2234
$wsExhibit.Range("F14").Value = 57
$wsExhibit.Range("F16").Value = 512
$wsExhibit.Range("F17").Value = 518
$wsExhibit.Range("F18").Value = 157
$wsExhibit.Range("F19").Value = 79
$wsExhibit.Range("F20").Value = 39
$wsExhibit.Range("F21").Value = 46
$wsExhibit.Range("F22").Value = 1719
$wsExhibit.Range("F23").Value = 3870
$wsExhibit.Range("F25").Value = 62
$wsExhibit.Range("F27").Value = 1153
$wsExhibit.Range("F28").Value = 216
$wsExhibit.Range("F29").Value = 2049
$wsExhibit.Range("F32").Value = 91
$wsExhibit.Range("F33").Value = 281
$wsExhibit.Range("F36").Value = 673
$wsExhibit.Range("F38").Value = 400

# Sheet "全部类型" (All Types) - update column F (想去人数 / "want to go" count)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 265
$wsAll.Range("F5").Value = 234
$wsAll.Range("F6").Value = 266
$wsAll.Range("F7").Value = 67
$wsAll.Range("F10").Value = 47
$wsAll.Range("F11").Value = 31
$wsAll.Range("F13").Value = 2234
$wsAll.Range("F14").Value = 57
$wsAll.Range("F17").Value = 512
$wsAll.Range("F18").Value = 518
$wsAll.Range("F19").Value = 157
$wsAll.Range("F20").Value = 79
$wsAll.Range("F21").Value = 39
$wsAll.Range("F22").Value = 46
$wsAll.Range("F23").Value = 1719
$wsAll.Range("F24").Value = 3870
$wsAll.Range("F26").Value = 62
$wsAll.Range("F28").Value = 1153
$wsAll.Range("F29").Value = 216
$wsAll.Range("F30").Value = 2049
$wsAll.Range("F33").Value = 91
$wsAll.Range("F34").Value = 281
$wsAll.Range("F37").Value = 673
$wsAll.Range("F39").Value = 401
